# Fixed Acceptance Test Plan
# For every in-scope "Sprint 2" row on the Test Plan sheet, mirror the
# Sprint-2 pass/fail result ("Pass") and comment ("TC 11/13/23") into the
# Sprint-3 columns (E = status, F = comment), matching what already exists
# further down the sheet (rows 54+).  Row 35 is a hidden row and is left
# untouched, and rows 13/14 are blank spacer rows with no data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

$rows = @(2,3,4,5,6,7,8,9,10,11,12,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "Pass"
    $ws.Cells.Item($r, 6).Value = "TC 11/13/23"
}

# Reproduce the author's final view state: scrolled back up a bit and with
# F31 as the active selection in the frozen right-hand pane.
$ws.Activate() | Out-Null
$ws.Range("F31").Select() | Out-Null
